# Atualização automática de pedidos - 30/05/2025 09:14
#
# 1) Fix RACK values that were previously stored as text to be real numbers
# 2) Append the new order "REQ-014" to both the "Pedidos" and "Itens" sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Pedidos"
# ---------------------------------------------------------------------------
$pedidos = $wb.Worksheets.Item("Pedidos")

# Existing row 14 (REQ-013) had its RACK value ("12") stored as text - store
# it as a real number instead.
$pedidos.Range("D14").Value = 12

# New row 15 - REQ-014
$pedidos.Range("A15").Value = "REQ-014"
$pedidos.Range("B15").Value = "30/05/2025 09:14"
$pedidos.Range("C15").Value = "Ford"
# RACK "10" keeps being stored as text (matches the source system behaviour)
$pedidos.Range("D15").Value = "'10"
$pedidos.Range("E15").Value = "R10-LA-A2"
$pedidos.Range("F15").Value = "washington vieira"
$pedidos.Range("H15").Value = "Pendente"

# ---------------------------------------------------------------------------
# Sheet "Itens"
# ---------------------------------------------------------------------------
$itens = $wb.Worksheets.Item("Itens")

# Existing row 14 (REQ-013) had its "seccao" value ("0.35") stored as text -
# store it as a real number instead.
$itens.Range("D14").Value = 0.35

# New row 15 - REQ-014
$itens.Range("A15").Value = "REQ-014"
$itens.Range("B15").Value = "LMPT2A-0.35-G-R"
$itens.Range("C15").Value = "180BA406965"
# "seccao" keeps being stored as text (matches the source system behaviour)
$itens.Range("D15").Value = "'0.35"
$itens.Range("E15").Value = "G-R"
$itens.Range("F15").Value = 1
